# Add season record columns (Wins, Losses, Ties) to the roster sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) ---
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Copy the header style from an existing header cell (A1) onto the new
# header cells so they match the bold/centered/bordered look of the rest
# of row 1.
$ws.Range("A1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# --- Data rows (rows 2-50): season record values ---
$lastRow = 50
$wins = $ws.Range("AD2:AD" + $lastRow)
$losses = $ws.Range("AE2:AE" + $lastRow)
$ties = $ws.Range("AF2:AF" + $lastRow)

$wins.Value = 74
$losses.Value = 88
$ties.Value = 0
